$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# This workbook stores the Price column (D) as text (e.g. "42.975.29", "98.31")
# even though the values look numeric. A plain Range.Value assignment would let
# Excel auto-convert numeric-looking strings into real numbers (losing trailing
# zeros such as "1.00" -> 1, or mis-parsing two-dot values). To faithfully keep
# these cells as text (matching the source data), force the cell to Text format
# before assigning the value, then clear the temporary formatting again so the
# cell keeps its original (default) style.

# Row 2
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '42.992.29'
$ws.Range("D2").ClearFormats()
$ws.Range("E2").Value = '  +0.57%  '

# Row 3
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '2.553.48'
$ws.Range("D3").ClearFormats()
$ws.Range("E3").Value = '  +0.74%  '

# Row 4
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '1.00'
$ws.Range("D4").ClearFormats()
$ws.Range("E4").Value = '  -0.14%  '

# Row 5
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '305.03'
$ws.Range("D5").ClearFormats()
$ws.Range("E5").Value = '  +2.31%  '

# Row 6
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '98.30'
$ws.Range("D6").ClearFormats()
$ws.Range("E6").Value = '  +6.67%  '

# Row 7
$ws.Range("E7").Value = '  +1.10%  '

# Row 8
$ws.Range("E8").Value = '  +0.04%  '

# Row 9
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.548'
$ws.Range("D9").ClearFormats()
$ws.Range("E9").Value = '  -0.09%  '

# Row 10
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '37.06'
$ws.Range("D10").ClearFormats()
$ws.Range("E10").Value = '  +3.79%  '

# Row 11
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.0828'
$ws.Range("D11").ClearFormats()
$ws.Range("E11").Value = '  +3.21%  '

# Row 12
$ws.Range("E12").Value = '  +1.84%  '

# Row 13
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '0.115'
$ws.Range("D13").ClearFormats()
$ws.Range("E13").Value = '  +2.41%  '

# Row 14
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '2.946.67'
$ws.Range("D14").ClearFormats()
$ws.Range("E14").Value = '  +0.67%  '

# Row 15
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '2.599.36'
$ws.Range("D15").ClearFormats()
$ws.Range("E15").Value = '  +3.24%  '

# Row 16
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '15.15'
$ws.Range("D16").ClearFormats()
$ws.Range("E16").Value = '  +7.85%  '

# Row 17
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '0.877'
$ws.Range("D17").ClearFormats()
$ws.Range("E17").Value = '  +1.07%  '

# Row 18
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '43.051.22'
$ws.Range("D18").ClearFormats()
$ws.Range("E18").Value = '  +0.52%  '

# Row 19
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '13.81'
$ws.Range("D19").ClearFormats()
$ws.Range("E19").Value = '  +7.54%  '

# Row 20
$ws.Range("E20").Value = '  +1.38%  '

# Row 21
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '6.59'
$ws.Range("D21").ClearFormats()
$ws.Range("E21").Value = '  +0.35%  '

# Row 22
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '71.99'
$ws.Range("D22").ClearFormats()
$ws.Range("E22").Value = '  +0.54%  '

# Row 23
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '255.27'
$ws.Range("D23").ClearFormats()
$ws.Range("E23").Value = '  -0.17%  '

# Row 24
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '2.97'
$ws.Range("D24").ClearFormats()
$ws.Range("E24").Value = '  +2.18%  '

# Row 25
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '2.09'
$ws.Range("D25").ClearFormats()
$ws.Range("E25").Value = '  -0.67%  '

# Row 26
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '28.12'
$ws.Range("D26").ClearFormats()
$ws.Range("E26").Value = '  -2.76%  '

# Row 27
$ws.Range("E27").Value = '  -0.05%  '

# Row 28
$ws.Range("E28").Value = '  +2.92%  '

# Row 29
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '37.89'
$ws.Range("D29").ClearFormats()
$ws.Range("E29").Value = '  +2.24%  '

# Row 30
$ws.Range("B30").Value = 'Filecoin'
$ws.Range("C30").Value = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '6.19'
$ws.Range("D30").ClearFormats()
$ws.Range("E30").Value = '  +4.84%  '

# Row 31
$ws.Range("B31").Value = 'Toncoin'
$ws.Range("C31").Value = 'https://coinranking.com/coin/67YlI0K1b+toncoin-ton'
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '2.08'
$ws.Range("D31").ClearFormats()
$ws.Range("E31").Value = '  -1.66%  '

# Row 32
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '158.61'
$ws.Range("D32").ClearFormats()
$ws.Range("E32").Value = '  +3.94%  '

# Row 33
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '19.64'
$ws.Range("D33").ClearFormats()
$ws.Range("E33").Value = '  +16.51%  '

# Row 34
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '2.16'
$ws.Range("D34").ClearFormats()
$ws.Range("E34").Value = '  +0.33%  '

# Row 35
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '0.0804'
$ws.Range("D35").ClearFormats()
$ws.Range("E35").Value = '  +1.09%  '

# Row 36
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '3.33'
$ws.Range("D36").ClearFormats()
$ws.Range("E36").Value = '  -1.64%  '

# Row 37
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '2.64'
$ws.Range("D37").ClearFormats()
$ws.Range("E37").Value = '  -3.81%  '

# Row 38
$ws.Range("E38").Value = '  +2.38%  '

# Row 39
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '25.66'
$ws.Range("D39").ClearFormats()
$ws.Range("E39").Value = '  +10.53%  '

# Row 40
$ws.Range("E40").Value = '  +0.11%  '

# Row 41
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '3.46'
$ws.Range("D41").ClearFormats()
$ws.Range("E41").Value = '  +1.49%  '

# Row 42
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '3.93'
$ws.Range("D42").ClearFormats()
$ws.Range("E42").Value = '  +1.18%  '

# Row 43
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '2.06'
$ws.Range("D43").ClearFormats()
$ws.Range("E43").Value = '  +29.43%  '

# Row 44
$ws.Range("B44").Value = 'VeChain'
$ws.Range("C44").Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '0.0307'
$ws.Range("D44").ClearFormats()
$ws.Range("E44").Value = '  -0.79%  '

# Row 45
$ws.Range("B45").Value = 'Maker'
$ws.Range("C45").Value = 'https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr'
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '2.093.25'
$ws.Range("D45").ClearFormats()
$ws.Range("E45").Value = '  +0.32%  '

# Row 46
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '1.00'
$ws.Range("D46").ClearFormats()
$ws.Range("E46").Value = '  +0.01%  '

# Row 47
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '86.92'
$ws.Range("D47").ClearFormats()
$ws.Range("E47").Value = '  +4.02%  '

# Row 48
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '9.00'
$ws.Range("D48").ClearFormats()
$ws.Range("E48").Value = '  +1.08%  '

# Row 49
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '75.43'
$ws.Range("D49").ClearFormats()
$ws.Range("E49").Value = '  +10.41%  '

# Row 50
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '2.803.75'
$ws.Range("D50").ClearFormats()
$ws.Range("E50").Value = '  +0.62%  '

# Row 51
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '0.193'
$ws.Range("D51").ClearFormats()
$ws.Range("E51").Value = '  +3.98%  '
